$d = $word.ActiveDocument

# --- Change 1: expand the "game will start" paragraph text ---
# Original: "The game will start in the park where you can walk around and talk with an NPC. ..."
# New:      "The game will start in a hotel room (the murder scene). You can exit the room and
#            walk around in the city. In the center of the city is a park where you can walk
#            around and talk with an NPC. ..."
$old1 = "The game will start in the park where you can walk around and talk with an NPC."
$new1 = "The game will start in a hotel room (the murder scene). You can exit the room and walk around in the city. In the center of the city is a park where you can walk around and talk with an NPC."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# Word keeps track of the last-edited location with the hidden "_GoBack" bookmark; after this
# edit it sits between "...In the center of the cit" and "y is a park...".
$marker = $d.Content
$marker.Find.Execute("In the center of the cit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackRange = $d.Range($marker.End, $marker.End)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- Change 2: add two new bullet items after "Walk through church door..." ---
$lastItem = $d.Content
$lastItem.Find.Execute("Walk through church door to get inside the church.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lastIndex = $lastItem.Paragraphs(1).Index

$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()
$pauseItem = $d.Paragraphs($lastIndex + 1)
$pauseItem.Range.Text = "Press esc to pause the game (from here you can quit or resume)"

$pauseItem.Range.InsertParagraphAfter()
$evidenceItem = $d.Paragraphs($lastIndex + 2)
$evidenceItem.Range.Text = "Press ‘e’ to see a list of all evidence you have collected so far"
